$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'51.788.23"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.10%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'2.803.97"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  +0.03%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'354.62"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -0.74%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'111.93"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +2.17%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'0.557"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +0.57%  "
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  -0.01%  "
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  +7.91%  "
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +1.63%  "
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -1.50%  "
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -0.86%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'20.03"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +2.32%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'7.77"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +2.03%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'3.243.96"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  +0.78%  "
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'2.805.60"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +1.22%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'0.945"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +1.16%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'51.807.53"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +0.24%  "
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  +1.16%  "
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  +3.41%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'13.69"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +3.73%  "
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +0.92%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'70.54"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  +0.39%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'268.77"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  +0.39%  "
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +1.19%  "
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'26.20"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -0.61%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'0.162"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -3.48%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'39.17"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  +11.55%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'10.40"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +2.04%  "
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +3.10%  "
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -0.74%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'52.27"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +0.71%  "
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +7.81%  "
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +6.03%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'0.0445"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  -0.72%  "
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.02%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'18.94"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +0.20%  "
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  +2.34%  "
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +0.41%  "
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  +1.34%  "
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -1.55%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'2.23"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +1.19%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'120.11"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +0.16%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'22.09"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  +1.03%  "
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  +4.57%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'2.120.01"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  +1.70%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'2.43"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  +6.62%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.959"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  +0.93%  "
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  -1.26%  "
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  +7.25%  "
$c.Style = "Normal"
